$d = $word.ActiveDocument

# 1. Remove the "Work-Authorization: EAD " run sequence from the header line,
#    leaving "| https://smakar20.github.io ..." etc. intact.
$d.Content.Find.Execute("Work-Authorization: EAD ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 2. Merge the split "Foundation program in Fi" / "nance (Oracle)." runs
#    (which used to be separated by the _GoBack bookmark) back into a
#    single run with the full text.
$d.Content.Find.Execute(" Foundation program in Finance (Oracle).", $true, $false, $false, $false, $false, $true, 1, $false, " Foundation program in Finance (Oracle).", 2) | Out-Null

# 3. Re-place the _GoBack bookmark so it spans the whole document body
#    (start of document through the end of the last paragraph), matching
#    where Word left the cursor for the final save. Adding a bookmark with
#    a name that already exists moves it, so this also removes the old
#    mid-document _GoBack (if any remained).
$full = $d.Range(0, $d.Content.End)
$d.Bookmarks.Add("_GoBack", $full) | Out-Null
